$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'65.429.51"
$ws.Range('E2').Value = '  -5.79%  '

$ws.Range('D3').Value = "'3.305.55"
$ws.Range('E3').Value = '  -6.55%  '

$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = "'555.50"
$ws.Range('E5').Value = '  -5.04%  '

$ws.Range('D6').Value = "'179.15"
$ws.Range('E6').Value = '  -8.44%  '

$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').Value = "'0.587"
$ws.Range('E8').Value = '  -3.74%  '

$ws.Range('D9').Value = "'3.302.81"
$ws.Range('E9').Value = '  -6.28%  '

$ws.Range('E10').Value = '  -10.85%  '

$ws.Range('D11').Value = "'0.581"
$ws.Range('E11').Value = '  -7.47%  '

$ws.Range('D12').Value = "'47.00"
$ws.Range('E12').Value = '  -10.75%  '

$ws.Range('D13').Value = "'0.0000262"
$ws.Range('E13').Value = '  -8.76%  '

$ws.Range('D14').Value = "'3.833.30"
$ws.Range('E14').Value = '  -6.56%  '

$ws.Range('E15').Value = '  -7.77%  '

$ws.Range('D16').Value = "'597.07"
$ws.Range('E16').Value = '  -10.02%  '

$ws.Range('D17').Value = "'17.98"
$ws.Range('E17').Value = '  -1.99%  '

$ws.Range('D18').Value = "'65.435.93"
$ws.Range('E18').Value = '  -5.88%  '

$ws.Range('D20').Value = "'3.295.48"
$ws.Range('E20').Value = '  -7.11%  '

$ws.Range('E21').Value = '  -8.85%  '

$ws.Range('D22').Value = "'0.896"
$ws.Range('E22').Value = '  -6.75%  '

$ws.Range('D23').Value = "'17.02"
$ws.Range('E23').Value = '  -5.71%  '

$ws.Range('D24').Value = "'101.73"
$ws.Range('E24').Value = '  -2.47%  '

$ws.Range('E25').Value = '  -6.62%  '

$ws.Range('D26').Value = "'3.96"
$ws.Range('E26').Value = '  -9.58%  '

$ws.Range('D27').Value = "'5.99"
$ws.Range('E27').Value = '  -0.38%  '

$ws.Range('E28').Value = '  -9.12%  '

$ws.Range('D29').Value = "'9.23"
$ws.Range('E29').Value = '  -8.61%  '

$ws.Range('D30').Value = "'8.61"
$ws.Range('E30').Value = '  -9.54%  '

$ws.Range('D31').Value = "'30.41"

$ws.Range('E32').Value = '  -11.85%  '

$ws.Range('D33').Value = "'6.19"
$ws.Range('E33').Value = '  -8.18%  '

$ws.Range('E34').Value = '  -6.50%  '

$ws.Range('E35').Value = '  -6.15%  '

$ws.Range('D36').Value = "'3.794.00"
$ws.Range('E36').Value = '  +0.45%  '

$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').Value = "'1.00"
$ws.Range('E37').Value = '  +0.06%  '

$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').Value = "'56.55"
$ws.Range('E38').Value = '  -8.47%  '

$ws.Range('D39').Value = "'518.62"
$ws.Range('E39').Value = '  +3.86%  '

$ws.Range('D40').Value = "'3.43"
$ws.Range('E40').Value = '  -9.22%  '

$ws.Range('D41').Value = "'0.0₃0706"
$ws.Range('E41').Value = '  -12.29%  '

$ws.Range('D42').Value = "'2.63"
$ws.Range('E42').Value = '  -8.78%  '

$ws.Range('B43').Value = 'CoreDAO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D43').Value = "'3.28"
$ws.Range('E43').Value = '  +23.08%  '

$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').Value = "'0.123"
$ws.Range('E44').Value = '  -8.23%  '

$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = "'31.82"
$ws.Range('E45').Value = '  -8.02%  '

$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').Value = "'0.337"
$ws.Range('E46').Value = '  -8.84%  '

$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = "'3.29"
$ws.Range('E47').Value = '  -3.54%  '

$ws.Range('D48').Value = "'0.0409"
$ws.Range('E48').Value = '  -8.56%  '

$ws.Range('E49').Value = '  -5.07%  '

$ws.Range('D50').Value = "'2.59"
$ws.Range('E50').Value = '  -9.86%  '

$ws.Range('D51').Value = "'0.998"
$ws.Range('E51').Value = '  -0.36%  '
